$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Allgemein")

# --- Row 22: was "42" (Oberschicht Veraenderungen...), becomes "2a" (OS Veraenderung (alle)) ---
$ws.Range("C22").Value = "2a"
$ws.Range("E22").Value = "OS Veränderung (alle)"
$ws.Range("H22").Value = "Differenz zwischen zwei TBk Bestandeskarten. Zeigt Abnahme und Zunahme der Oberschicht für alle Bestände."
$ws.Range("I22").Clear()
$ws.Range("L22").Clear()

# --- Row 23: was "42" (OS Veraenderung (alle)), becomes "2b" (Oberschicht Veraenderungen...) ---
$ws.Range("C23").Value = "2b"
$ws.Range("E23").Value = "Oberschicht Veränderungen (hdom>30m) [OS VÄ]"
$ws.Range("G23").Formula = '=CONCAT(B23,C23,"_",IF(D23="","",CONCAT(D23," ")),IF(F23="","",CONCAT(F23," ")),E23)'
$ws.Range("H23").Value = "Differenz zwischen zwei TBk Bestandeskarten. Zeigt Abnahme und Zunahme der Oberschicht und geräumte Bestände (hdom > 30m)"
$ws.Range("I23").Value = "x"
$ws.Range("L23").Value = "x"

# --- Column G slightly wider ---
$ws.Columns.Item(7).ColumnWidth = 33.62

# --- Cursor/selection moved to G9 ---
$ws.Activate()
$ws.Range("G9").Select()
